$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean
$ws.Cells.Clear()

# New header row: Id, ParentId, Name, Options
$ws.Cells.Item(1, 1).Value = "Id"
$ws.Cells.Item(1, 2).Value = "ParentId"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Options"

# New staff-category rows (Id, ParentId, Name) - Options left blank
$data = @(
    @("1", "NULL", "Администрация"),
    @("2", "NULL", "Врачебный персонал"),
    @("3", "NULL", "Средний медицинский персонал"),
    @("4", "NULL", "Младший медицинский персонал"),
    @("5", "NULL", "Технический персонал")
)

$r = 2
foreach ($row in $data) {
    # Leading apostrophe forces the numeric-looking Id to be stored as text,
    # matching the original workbook's convention of text Id values.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Drop the implicit "Text" number format the apostrophe-prefixed entries
# picked up so every cell keeps the workbook's default style.
$ws.Cells.ClearFormats()
